$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

# "Förändrad" (column C) date was updated from 2023-09-06 (45175) to
# 2023-09-14 (45183) for every data row (rows 2-12).
for ($row = 2; $row -le 12; $row++) {
    $ws.Cells.Item($row, 3).Value = 45183
}
